{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Find index of first target paragraph, then also delete the blank paragraph right before it.\nlet idx = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === targets[0]) {\n    idx = i;\n    break;\n  }\n}\n\nif (idx >= 0) {\n  // delete the blank paragraph immediately preceding (if blank)\n  if (idx - 1 >= 0 && paras.items[idx - 1].text.trim() === \"\") {\n    paras.items[idx - 1].delete();\n  }\n  paras.items[idx].delete();\n  if (idx + 1 < paras.items.length && paras.items[idx + 1].text === targets[1]) {\n    paras.items[idx + 1].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$target1 = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$target2 = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$paras = @($d.Paragraphs)\n\n# Locate the index (0-based) of the \"Ver no Jupiter...\" paragraph.\n$idx = -1\nfor ($i = 0; $i -lt $paras.Count; $i++) {\n    $text = $paras[$i].Range.Text.TrimEnd([char]0x0D, [char]0x07)\n    if ($text -eq $target1) {\n        $idx = $i\n        break\n    }\n}\n\nif ($idx -ge 0) {\n    # Delete from the bottom up so earlier indices/ranges stay valid.\n    # 1) \"Ver no Jupiter...\" line (idx)\n    # 2) the blank paragraph right before it, if blank (idx - 1)\n    # 3) the following copyright line, if it matches (idx + 1)\n    if (($idx + 1) -lt $paras.Count) {\n        $nextText = $paras[$idx + 1].Range.Text.TrimEnd([char]0x0D, [char]0x07)\n        if ($nextText -eq $target2) {\n            $paras[$idx + 1].Range.Delete()\n        }\n    }\n\n    $paras[$idx].Range.Delete()\n\n    if ($idx -gt 0) {\n        $prevText = $paras[$idx - 1].Range.Text.TrimEnd([char]0x0D, [char]0x07)\n        if ($prevText.Trim() -eq \"\") {\n            $paras[$idx - 1].Range.Delete()\n        }\n    }\n}\n"}
